# Apply edits described by the diff:
# 1. On the "Config" sheet, insert a new row above row 5 containing
#    "bus_protocol" / "apb", shifting num_write_ports..byte_enable down by one.
# 2. On the "Registers" sheet, change C4 from "ReadClean" to "Write1Clear".

$wb = $excel.ActiveWorkbook

$configSheet = $wb.Worksheets.Item("Config")
$configSheet.Rows.Item(5).Insert()
$configSheet.Cells.Item(5, 1).Value = "bus_protocol"
$configSheet.Cells.Item(5, 2).Value = "apb"

$registersSheet = $wb.Worksheets.Item("Registers")
$registersSheet.Cells.Item(4, 3).Value = "Write1Clear"
